$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the embedded SQL queries (B2:B7 and C2) so the join conditions
#    use the renamed id columns (study_id / participant_id) instead of the
#    generic "id" column.
# ---------------------------------------------------------------------------
$replacements = @(
    @{ Old = 'std\.id = prt\."study\.id"';               New = 'std.study_id = prt."study.study_id"' },
    @{ Old = 'prt\.id = dgn\."participant\.id"';          New = 'prt.participant_id = dgn."participant.participant_id"' },
    @{ Old = 'prt\.id = trt\."participant\.id"';          New = 'prt.participant_id = trt."participant.participant_id"' },
    @{ Old = 'prt\.id = trr\."participant\.id"';          New = 'prt.participant_id = trr."participant.participant_id"' },
    @{ Old = 'prt\.id = srv\."participant\.id"';          New = 'prt.participant_id = srv."participant.participant_id"' },
    @{ Old = 'std\.id = rfs\."study\.id"';                New = 'std.study_id = rfs."study.study_id"' }
)

$queryCells = @("B2", "B3", "B4", "B5", "B6", "B7", "C2")

foreach ($addr in $queryCells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null) {
        foreach ($r in $replacements) {
            $text = $text -replace $r.Old, $r.New
        }
        $cell.Value2 = $text
    }
}

# ---------------------------------------------------------------------------
# 2. Resize column C (drop the auto "best fit" width, switch to an explicit
#    wider width ~67.16 characters) and refresh the view: scroll back to the
#    top and move the selection from C7 to B2.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 66.3

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
